# Remove the "culture_collection" attribute column (column AH) from the
# MIGS sediment sample sheet, re-derived from the 2017 INSDC review ("culture
# collection を MIxS から再度削除"). Column AH currently holds the header
# "culture_collection" (row 15) together with its cell comment describing the
# field; every column from AI through CX must shift one position to the left
# to take its place, and the vacated last column (CX) is dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = 15
$firstShiftCol = 35   # column AI  (first column to slide left)
$lastShiftCol = 102   # column CX  (last column, becomes empty afterwards)

# 1) Slide every field-description comment one column to the left (AI15's
#    comment becomes AH15's comment, AJ15's becomes AI15's, and so on) so the
#    comments stay attached to the same field names as the header text moves.
for ($col = $firstShiftCol; $col -le $lastShiftCol; $col++) {
    $srcComment = $ws.Cells.Item($headerRow, $col).Comment
    $text = $srcComment.Text()
    $dstComment = $ws.Cells.Item($headerRow, $col - 1).Comment
    $dstComment.Text($text) | Out-Null
}

# The old last-column comment (CX15) is now duplicated onto CW15, so drop it.
$ws.Cells.Item($headerRow, $lastShiftCol).Comment.Delete()

# 2) Slide the header labels themselves one column to the left to match.
for ($col = $firstShiftCol; $col -le $lastShiftCol; $col++) {
    $srcVal = $ws.Cells.Item($headerRow, $col).Value()
    $ws.Cells.Item($headerRow, $col - 1).Value = $srcVal
}

# The last column (formerly CX) is now a duplicate of CW; clear it out
# entirely so the row/sheet no longer references that extra column.
$ws.Cells.Item($headerRow, $lastShiftCol).Clear()
